{"js": "// --- Change 1: split the run right before \"herefore, we\" and relocate\n//     the auto \"_GoBack\" bookmark to that point (removing it from wherever\n//     it used to sit -- a bookmark name is unique per document).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst splitResults = context.document.body.search(\"current value of -10.  T\", { matchCase: true });\nsplitResults.load(\"items\");\nawait context.sync();\n\nif (splitResults.items.length > 0) {\n  const splitPoint = splitResults.items[0].getRange(\"End\");\n  splitPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- Change 2: the bulleted/highlighted question paragraph is cleared\n//     out entirely -- no more numbering, no more highlighted run, just\n//     an empty paragraph that keeps the spacing and picks up a manual\n//     left indent instead of the list's indent.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetText = \"Can you think of a real domain in which step costs are such as to cause looping?\";\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    const replacement = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:pPr><w:spacing w:before=\"100\" w:beforeAutospacing=\"1\" w:after=\"100\" w:afterAutospacing=\"1\"/><w:ind w:left=\"720\"/><w:rPr><w:rFonts w:ascii=\"Times\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times\" w:cs=\"Times New Roman\"/><w:color w:val=\"000000\"/><w:sz w:val=\"27\"/><w:szCs w:val=\"27\"/></w:rPr></w:pPr></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n    paragraphs.items[i].getRange(\"Whole\").insertOoxml(replacement, Word.InsertLocation.replace);\n    await context.sync();\n    break;\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: split the run right before \"herefore, we\" and drop the\n#     auto \"_GoBack\" bookmark there (this also removes it from wherever\n#     it used to sit, since a bookmark name is unique per document).\n$r = $d.Content\n$found = $r.Find.Execute(\"current value of -10.  T\")\nif ($found) {\n    $r.Collapse(0)  # wdCollapseEnd\n    $null = $d.Bookmarks.Add(\"_GoBack\", $r)\n}\n\n# --- Change 2: the bulleted/highlighted question paragraph is cleared\n#     out entirely -- no more numbering, no more highlighted run, just\n#     an empty paragraph that keeps the spacing and picks up a manual\n#     left indent instead of the list's indent.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Can you think of a real domain in which step costs are such as to cause looping?*\") {\n        $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:pPr><w:spacing w:before=\"100\" w:beforeAutospacing=\"1\" w:after=\"100\" w:afterAutospacing=\"1\"/><w:ind w:left=\"720\"/><w:rPr><w:rFonts w:ascii=\"Times\" w:eastAsia=\"Times New Roman\" w:hAnsi=\"Times\" w:cs=\"Times New Roman\"/><w:color w:val=\"000000\"/><w:sz w:val=\"27\"/><w:szCs w:val=\"27\"/></w:rPr></w:pPr></w:p>'\n        $null = $p.Range.InsertXML($xml)\n        break\n    }\n}\n"}
